$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (78) duplicating the last existing data row (77),
# as new sensor data pulled from Adafruit IO.
$row = 78

$ws.Cells.Item($row, 3).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"
$ws.Cells.Item($row, 3).Value = "25"
$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"

# Restore the default (General) style on C78 now that the text value is
# locked in, so the new row carries no style override - matching the
# other plain data rows.
$ws.Cells.Item($row, 3).ClearFormats()
